$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the survey question text: the line break moved from after "have" to after "countries"
$ws.Range("A2").Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""

# Re-fit the row height so no stray custom row height is introduced by the wrap
$ws.Rows(2).AutoFit()

# Update the recalculated mean values
$ws.Range("B2").Value = 0.354477546162626
$ws.Range("K2").Value = 0.330419625213169
$ws.Range("L2").Value = 0.55583349037107
$ws.Range("N2").Value = 0.105209195768849
